# Applies the "Added login signup requirements" edit: a new LS_01..LS_16
# requirements block (rows 91-107) appended to the Login/Signup section.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write column-by-column (C, then A, then B, then D) so new shared-string
# entries are interned in the same order the original workbook used for
# every other section of this sheet (matches the source diff ordering).

# --- Column C: requirement descriptions ---
$ws.Range("C91").Value2 = "System shall have a username box for the returning user to enter thier unique identifier"
$ws.Range("C92").Value2 = "System shall have a password box for the returning user to enter their created password."
$ws.Range("C93").Value2 = "System shall have a login button for directing returning users to login page."
$ws.Range("C94").Value2 = "System shall have a sign up button for directing new users to new account page."
$ws.Range("C95").Value2 = "System shall have a game title (logo)"
$ws.Range("C96").Value2 = "System shall have a sign up information input to create new accounts for new users."
$ws.Range("C97").Value2 = "System users shall have a unique username with at least 4 characters."
$ws.Range("C98").Value2 = "System users shall have a password with at least 6 unique characters."
$ws.Range("C99").Value2 = "System shall display a message if no account is found matching the entered username."
$ws.Range("C100").Value2 = "System will comunicate with Page Manager to send data to database"
$ws.Range("C101").Value2 = "System shall have a username box for the new user to enter thier unique identifier"
$ws.Range("C102").Value2 = "System shall have a password box for the new user to enter their created password."
$ws.Range("C103").Value2 = "System shall have a functioning button to log the user in using their entered credentials."
$ws.Range("C104").Value2 = "System shall have a functioning button to create a user account using their entered credentials."
$ws.Range("C105").Value2 = "System shall have a function to show the signup page and simultaneuously hide the login screen upon click. Function for Requirement(4)"
$ws.Range("C106").Value2 = "System shall have a function to show the login page and simultaneuously hide the signup screen upon click. Function for Requirement(3)"

# --- Column A: requirement ids ---
$ws.Range("A91").Value2 = "LS_01"
$ws.Range("A92").Value2 = "LS_02"
$ws.Range("A93").Value2 = "LS_03"
$ws.Range("A94").Value2 = "LS_04"
$ws.Range("A95").Value2 = "LS_05"
$ws.Range("A96").Value2 = "LS_06"
$ws.Range("A97").Value2 = "LS_07"
$ws.Range("A98").Value2 = "LS_08"
$ws.Range("A99").Value2 = "LS_09"
$ws.Range("A100").Value2 = "LS_10"
$ws.Range("A101").Value2 = "LS_11"
$ws.Range("A102").Value2 = "LS_12"
$ws.Range("A103").Value2 = "LS_13"
$ws.Range("A104").Value2 = "LS_14"
$ws.Range("A105").Value2 = "LS_15"
$ws.Range("A106").Value2 = "LS_16"

# --- Column B: section labels ---
$ws.Range("B91").Value2 = "01 - Login/Signup"
$ws.Range("B92").Value2 = "01 - Login/Signup"
$ws.Range("B93").Value2 = "1 - Login/Signup"
$ws.Range("B94").Value2 = "1 - Login/Signup"
$ws.Range("B95").Value2 = "1 - Login/Signup"
$ws.Range("B96").Value2 = "1 - Login/Signup"
$ws.Range("B97").Value2 = "1 - Login/Signup"
$ws.Range("B98").Value2 = "1 - Login/Signup"
$ws.Range("B99").Value2 = "1 - Login/Signup"
$ws.Range("B100").Value2 = "1 - Login/Signup"
$ws.Range("B101").Value2 = "1 - Login/Signup"
$ws.Range("B102").Value2 = "1 - Login/Signup"
$ws.Range("B103").Value2 = "1 - Login/Signup"
$ws.Range("B104").Value2 = "1 - Login/Signup"
$ws.Range("B105").Value2 = "1 - Login/Signup"
$ws.Range("B106").Value2 = "1 - Login/Signup"

# --- Column D: short cross-reference tags ---
$ws.Range("D91").Value2 = "login"
$ws.Range("D92").Value2 = "login"
$ws.Range("D93").Value2 = "btn"
$ws.Range("D94").Value2 = "btn"
$ws.Range("D95").Value2 = "img"
$ws.Range("D96").Value2 = "new_account"
$ws.Range("D97").Value2 = "uniqueUser, span"
$ws.Range("D98").Value2 = "uniquePass, span"
$ws.Range("D99").Value2 = "span"
$ws.Range("D100").Value2 = "sendInfo"
$ws.Range("D101").Value2 = "new_account"
$ws.Range("D102").Value2 = "new_account"
$ws.Range("D103").Value2 = "login (id)"
$ws.Range("D104").Value2 = "new_account (id)"
$ws.Range("D105").Value2 = "auth.js"
$ws.Range("D106").Value2 = "auth.js"

# --- Formatting ---

# Column D (rows 91-106): centered horizontally + vertically (new style).
$dRange = $ws.Range("D91:D106")
$dRange.Style = "Normal"
$dRange.HorizontalAlignment = -4108   # xlCenter
$dRange.VerticalAlignment = -4108     # xlCenter

# C105/C106 wrap (two-line requirement text) + matching row height.
$wrapRange = $ws.Range("C105:C106")
$wrapRange.Style = "Normal"
$wrapRange.WrapText = $true
$ws.Rows.Item(105).RowHeight = 28.8
$ws.Rows.Item(106).RowHeight = 28.8

# Trailing blank styled cell that closes the block (matches source row 107).
$ws.Range("C107").Style = "Normal"
$ws.Range("C107").WrapText = $true

# Restore the view: scrolled down to the new block, E92 selected.
$ws.Range("E92").Select()
